$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) updates that look like plain numbers stay as
# literal text (matching the original inlineStr formatting, e.g. "246.98")
# instead of being auto-converted to floating point numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.902.49"
$ws.Range("E2").Value = "  -0.75%  "
$ws.Range("D3").Value = "1.741.26"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "246.98"
$ws.Range("E5").Value = "  +4.85%  "
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("E7").Value = "  -4.50%  "
$ws.Range("D8").Value = "0.2717"
$ws.Range("E8").Value = "  -2.97%  "
$ws.Range("D9").Value = "0.06178"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").Value = "1.741.61"
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("D11").Value = "0.07231"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "15.07"
$ws.Range("E12").Value = "  -2.69%  "
$ws.Range("D13").Value = "0.6459"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").Value = "4.623"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "77.46"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "25.912.61"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "11.81"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").Value = "0.000006798"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("D21").Value = "1.967.38"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "4.274"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").Value = "8.622"
$ws.Range("E23").Value = "  -1.39%  "
$ws.Range("D24").Value = "5.376"
$ws.Range("D25").Value = "136.08"
$ws.Range("E25").Value = "  -1.85%  "
$ws.Range("D26").Value = "1.501"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "15.23"
$ws.Range("E27").Value = "  -0.57%  "
$ws.Range("D28").Value = "1.763"
$ws.Range("E28").Value = "  -2.81%  "
$ws.Range("D29").Value = "105.57"
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("D30").Value = "3.909"
$ws.Range("E30").Value = "  +2.27%  "
$ws.Range("D31").Value = "0.08221"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "3.626"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").Value = "0.04668"
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").Value = "2.657"
$ws.Range("E34").Value = "  +0.42%  "
$ws.Range("D35").Value = "0.9919"
$ws.Range("E35").Value = "  -1.81%  "
$ws.Range("D36").Value = "0.6197"
$ws.Range("E36").Value = "  -2.74%  "
$ws.Range("D37").Value = "2.730"
$ws.Range("E37").Value = "  +0.66%  "
$ws.Range("D38").Value = "0.01599"
$ws.Range("E38").Value = "  -0.39%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "1.912"
$ws.Range("E39").Value = "  -2.44%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "1.001"
$ws.Range("E40").Value = "  +0.20%  "
$ws.Range("D41").Value = "98.88"
$ws.Range("E41").Value = "  -1.79%  "
$ws.Range("D42").Value = "0.7571"
$ws.Range("D43").Value = "0.3857"
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("D44").Value = "4.979"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("D45").Value = "0.1134"
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("D46").Value = "6.253"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("E47").Value = "  +1.80%  "
$ws.Range("D48").Value = "0.05233"
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("E49").Value = "  -2.01%  "
$ws.Range("D50").Value = "7.485"
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("D51").Value = "0.3406"
$ws.Range("E51").Value = "  -1.64%  "
